$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "第一次修改"
$ws.Range("A4").Value = "第三次修改"
$ws.Range("B4").Value = "第二次的内容丢失了"

$ws.Range("B4").Select()
